$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: add A12 "decimal" (B12 "dry_bag_wt_grams_000" already present)
$ws.Range("A12").Value = "decimal"

# Row 13 (new): date / decomp_bag_collect_date / label
$ws.Range("A13").Value = "date"
$ws.Range("B13").Value = "decomp_bag_collect_date"
$ws.Range("C13").Value = "Date the bags in the form were collected"

# Row 14: add B14 note (A14 already present)
$ws.Range("B14").Value = "ex; barcode_bag and barcode_yield"

# Update the selection to match the new active cell
$ws.Range("B15").Select()
